$d = $word.ActiveDocument

$d.Content.Find.Execute("390÷3=130, 0", $true, $false, $false, $false, $false, $true, 1, $false, "363÷8=45, 3", 2) | Out-Null
$d.Content.Find.Execute("625÷3=208, 1", $true, $false, $false, $false, $false, $true, 1, $false, "268÷5=53, 3", 2) | Out-Null
$d.Content.Find.Execute("925÷4=231, 1", $true, $false, $false, $false, $false, $true, 1, $false, "350÷2=175, 0", 2) | Out-Null
$d.Content.Find.Execute("143÷2=71, 1", $true, $false, $false, $false, $false, $true, 1, $false, "458÷3=152, 2", 2) | Out-Null
$d.Content.Find.Execute("509÷3=169, 2", $true, $false, $false, $false, $false, $true, 1, $false, "420÷8=52, 4", 2) | Out-Null
$d.Content.Find.Execute("741÷3=247, 0", $true, $false, $false, $false, $false, $true, 1, $false, "325÷2=162, 1", 2) | Out-Null
$d.Content.Find.Execute("724÷7=103, 3", $true, $false, $false, $false, $false, $true, 1, $false, "710÷4=177, 2", 2) | Out-Null
$d.Content.Find.Execute("885÷7=126, 3", $true, $false, $false, $false, $false, $true, 1, $false, "794÷3=264, 2", 2) | Out-Null
$d.Content.Find.Execute("925÷7=132, 1", $true, $false, $false, $false, $false, $true, 1, $false, "836÷4=209, 0", 2) | Out-Null
$d.Content.Find.Execute("823÷5=164, 3", $true, $false, $false, $false, $false, $true, 1, $false, "717÷7=102, 3", 2) | Out-Null
$d.Content.Find.Execute("371÷9=41, 2", $true, $false, $false, $false, $false, $true, 1, $false, "248÷3=82, 2", 2) | Out-Null
$d.Content.Find.Execute("471÷9=52, 3", $true, $false, $false, $false, $false, $true, 1, $false, "607÷3=202, 1", 2) | Out-Null
$d.Content.Find.Execute("139÷2=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "669÷4=167, 1", 2) | Out-Null
$d.Content.Find.Execute("484÷3=161, 1", $true, $false, $false, $false, $false, $true, 1, $false, "710÷8=88, 6", 2) | Out-Null
$d.Content.Find.Execute("228÷9=25, 3", $true, $false, $false, $false, $false, $true, 1, $false, "993÷7=141, 6", 2) | Out-Null
$d.Content.Find.Execute("796÷5=159, 1", $true, $false, $false, $false, $false, $true, 1, $false, "557÷5=111, 2", 2) | Out-Null
$d.Content.Find.Execute("893÷3=297, 2", $true, $false, $false, $false, $false, $true, 1, $false, "481÷2=240, 1", 2) | Out-Null
$d.Content.Find.Execute("123÷4=30, 3", $true, $false, $false, $false, $false, $true, 1, $false, "699÷8=87, 3", 2) | Out-Null
$d.Content.Find.Execute("743÷6=123, 5", $true, $false, $false, $false, $false, $true, 1, $false, "506÷5=101, 1", 2) | Out-Null
$d.Content.Find.Execute("398÷2=199, 0", $true, $false, $false, $false, $false, $true, 1, $false, "330÷2=165, 0", 2) | Out-Null
$d.Content.Find.Execute("106÷2=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "867÷5=173, 2", 2) | Out-Null
$d.Content.Find.Execute("562÷9=62, 4", $true, $false, $false, $false, $false, $true, 1, $false, "666÷2=333, 0", 2) | Out-Null
$d.Content.Find.Execute("624÷2=312, 0", $true, $false, $false, $false, $false, $true, 1, $false, "823÷9=91, 4", 2) | Out-Null
$d.Content.Find.Execute("367÷5=73, 2", $true, $false, $false, $false, $false, $true, 1, $false, "883÷7=126, 1", 2) | Out-Null
$d.Content.Find.Execute("497÷8=62, 1", $true, $false, $false, $false, $false, $true, 1, $false, "849÷3=283, 0", 2) | Out-Null
